# Apply "Minor updates, added some data to Habitat_Raw" edits.
# Updates a handful of rating cells on the active sheet and clears a few
# cells that are no longer populated for certain reaches.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 281 (Libby Creek 01) ---
$ws.Range("C281").Value = "At Risk"
$ws.Range("D281").Value = "At Risk"

# --- Row 282 (Libby Creek 02) ---
$ws.Range("C282").Value = "At Risk"
$ws.Range("D282").Value = "At Risk"

# --- Row 547 (Twisp River Upper 04) ---
$ws.Range("E547").Value = "Unacceptable"
$ws.Range("I547").Value = "At Risk"

# --- Row 601 (Wolf Creek 01) ---
$ws.Range("E601").Value = "Unacceptable"
$ws.Range("G601").Value = "Unacceptable"
$ws.Range("H601").Value = "Unacceptable"
$ws.Range("I601").Value = "Unacceptable"
$ws.Range("J601").Value = "Unacceptable"
$ws.Range("K601").Value = "Unacceptable"
$ws.Range("L601").Value = "At Risk"
$ws.Range("M601").Value = "Unacceptable"
$ws.Range("N601").Value = "At Risk"
$ws.Range("R601").ClearContents()
$ws.Range("AC601").ClearContents()

# --- Row 602 (Wolf Creek 02) ---
$ws.Range("F602").Value = "Unacceptable"
$ws.Range("G602").Value = "Adequate"
$ws.Range("J602").Value = "Adequate"
$ws.Range("K602").Value = "Adequate"
$ws.Range("L602").Value = "Adequate"
$ws.Range("M602").Value = "At Risk"
$ws.Range("N602").Value = "At Risk"
$ws.Range("AC602").ClearContents()

# --- Row 603 (Wolf Creek 03) ---
$ws.Range("E603").Value = "Adequate"
$ws.Range("F603").Value = "Unacceptable"
$ws.Range("G603").Value = "Adequate"
$ws.Range("H603").Value = "Adequate"
$ws.Range("I603").Value = "Adequate"
$ws.Range("J603").Value = "Adequate"
$ws.Range("K603").Value = "Adequate"
$ws.Range("L603").Value = "Adequate"
$ws.Range("M603").Value = "Adequate"
$ws.Range("N603").Value = "Adequate"
$ws.Range("R603").ClearContents()
$ws.Range("AC603").ClearContents()
